# Fruta / hortaliza, semanal
# Swap data between the two groups of rows: {2, 10, 4} and {3, 12, 5}
# Each row's D, M, N, O, P, R, S values are replaced with the values
# currently held by the "source" row, forming two 3-cycles:
#   2 <- 10 <- 4 <- 2   and   3 <- 12 <- 5 <- 3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move together as a row's "data" for this swap
$cols = @("D", "M", "N", "O", "P", "R", "S")

function Get-RowData($ws, $row, $cols) {
    $data = @{}
    foreach ($col in $cols) {
        $data[$col] = $ws.Range("$col$row").Value2
    }
    return $data
}

# Snapshot original values of every row involved before overwriting any of them
$orig2  = Get-RowData $ws 2  $cols
$orig3  = Get-RowData $ws 3  $cols
$orig4  = Get-RowData $ws 4  $cols
$orig5  = Get-RowData $ws 5  $cols
$orig10 = Get-RowData $ws 10 $cols
$orig12 = Get-RowData $ws 12 $cols

function Set-RowData($ws, $row, $cols, $data) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value2 = $data[$col]
    }
}

# row 2  <- old row 10
Set-RowData $ws 2  $cols $orig10
# row 10 <- old row 4
Set-RowData $ws 10 $cols $orig4
# row 4  <- old row 2
Set-RowData $ws 4  $cols $orig2

# row 3  <- old row 12
Set-RowData $ws 3  $cols $orig12
# row 12 <- old row 5
Set-RowData $ws 12 $cols $orig5
# row 5  <- old row 3
Set-RowData $ws 5  $cols $orig3
